$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Test case "test_CreateUser" (row 28) passed: update Outcome (F28) to match the
# Expected Outcome (E28), and record the Justification (G28).
$ws.Range("F28").Value = "A new user is created"
$ws.Range("G28").Value = "Based on the given source code, this function has already been implemented"

# Reflect the reviewer's final scroll position / selection on the sheet.
$ws.Range("D17").Select() | Out-Null
$ws.Range("G27").Select() | Out-Null
